# Katrevula_LabExam03Grading.xlsx - grading pass: fill in "Points for grading"
# (column E) for the first two rubric sections to mirror the max points
# already recorded in column D ("Total Points"), i.e. full marks awarded.
# Corresponding "Total" formulas (E7, E15, E38) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section 1 (Generic rubric, rows 3-6)
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Section 2 (Customer Class rubric, rows 10-14)
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Leave the cursor/selection where the grader last worked
$ws.Range("E15").Select()
